$d = $word.ActiveDocument

# Locate the paragraph that currently reads exactly "Added ability to stack items".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Added ability to stack items`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph 'Added ability to stack items'"
}

# Pull the paragraph's own <w:p ...> opening tag (with its original identity
# attributes: w14:paraId, w14:textId, w:rsidR, w:rsidP, ...) so the rewritten
# first paragraph keeps them instead of getting a bare <w:p>.
$owx = $target.Range.WordOpenXML
$openTag = "<w:p>"
if ($owx -match '(<w:p[ >][^>]*>)') {
    $openTag = $matches[1]
}

# Grab the paragraph's existing <w:pPr> (list style + numbering) so the
# brand-new sibling paragraph we add below uses the same list formatting.
$pPr = ""
if ($owx -match '(<w:pPr>.*?</w:pPr>)') {
    $pPr = $matches[1]
}

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>$openTag$pPr<w:r><w:t xml:space="preserve">Added ability to stack </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>items</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p>$pPr<w:r><w:t xml:space="preserve">Added ability to use </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>items</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

[void]$target.Range.InsertXML($xml)
